$wb = $excel.ActiveWorkbook

# --- Fix duplicate "id" values in the anomaly example of "employees" (Blad1) ---
$ws1 = $wb.Worksheets.Item("Blad1")

# First employees table (Table1, B5:E12) - id column B, rows 7-12
$ws1.Range("B7").Value = 1
$ws1.Range("B8").Value = 2
$ws1.Range("B9").Value = 3
$ws1.Range("B10").Value = 4
$ws1.Range("B11").Value = 5
$ws1.Range("B12").Value = 6

# Second employees table (Table14, B20:E27) - id column B, rows 22-27
$ws1.Range("B22").Value = 1
$ws1.Range("B23").Value = 2
$ws1.Range("B24").Value = 3
$ws1.Range("B25").Value = 4
$ws1.Range("B26").Value = 5
$ws1.Range("B27").Value = 6

# --- Restore the active sheet / selections to match the saved view state ---
$ws2 = $wb.Worksheets.Item("From 0NF to 3NF")
[void]$ws2.Range("N5").Select()

[void]$ws1.Activate()
[void]$ws1.Range("F16").Select()
